$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ2 = $wb.Worksheets.Item(2)   # original "2022-Q2" sheet; stays untouched

# 1. Duplicate the "2022-Q2" sheet (keeps header/column styling) and place it
#    right after "总计"; this becomes the new "2022-Q3" sheet.
$sheetQ2.Copy($null, $sheetTotal)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Extend the column-A index style down to the extra rows (11-16) the new
# dataset needs beyond the original 10-row range.
$newSheet.Range("A10:H10").Copy()
$newSheet.Range("A11:H16").PasteSpecial(-4122)

# 2. Overwrite header row (identical text, kept for clarity/robustness)
$newSheet.Range("B1").Value = '基金代码'
$newSheet.Range("C1").Value = '基金名称'
$newSheet.Range("D1").Value = '基金规模'
$newSheet.Range("E1").Value = '股票总仓位'
$newSheet.Range("F1").Value = '仓位占比'
$newSheet.Range("G1").Value = '持有市值(亿元)'
$newSheet.Range("H1").Value = '仓位排名'

# 3. Write the 2022-Q3 fund rows (fund code + D/E/F/G as text via quote-prefix,
#    matching the source formatting; H is numeric; G16 is the lone numeric-0 exception)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = '''005775'
$newSheet.Range("C2").Value = '中加转型动力灵活配置混合A'
$newSheet.Range("D2").Value = '''5.34'
$newSheet.Range("E2").Value = '''50.55'
$newSheet.Range("F2").Value = '''2.87'
$newSheet.Range("G2").Value = '''0.1533'
$newSheet.Range("H2").Value = 5
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = '''005561'
$newSheet.Range("C3").Value = '创金合信中证红利低波动指数A'
$newSheet.Range("D3").Value = '''3.32'
$newSheet.Range("E3").Value = '''94.12'
$newSheet.Range("F3").Value = '''2.40'
$newSheet.Range("G3").Value = '''0.0797'
$newSheet.Range("H3").Value = 5
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = '''512890'
$newSheet.Range("C4").Value = '华泰柏瑞中证红利低波动ETF'
$newSheet.Range("D4").Value = '''2.60'
$newSheet.Range("E4").Value = '''99.50'
$newSheet.Range("F4").Value = '''2.57'
$newSheet.Range("G4").Value = '''0.0668'
$newSheet.Range("H4").Value = 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = '''005562'
$newSheet.Range("C5").Value = '创金合信中证红利低波动指数C'
$newSheet.Range("D5").Value = '''2.19'
$newSheet.Range("E5").Value = '''94.12'
$newSheet.Range("F5").Value = '''2.40'
$newSheet.Range("G5").Value = '''0.0526'
$newSheet.Range("H5").Value = 5
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = '''009927'
$newSheet.Range("C6").Value = '工银瑞信聚利18个月定期开放混合A'
$newSheet.Range("D6").Value = '''2.16'
$newSheet.Range("E6").Value = '''26.72'
$newSheet.Range("F6").Value = '''1.94'
$newSheet.Range("G6").Value = '''0.0419'
$newSheet.Range("H6").Value = 3
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = '''012072'
$newSheet.Range("C7").Value = '中加喜利回报一年持有期混合C'
$newSheet.Range("D7").Value = '''2.21'
$newSheet.Range("E7").Value = '''38.64'
$newSheet.Range("F7").Value = '''1.71'
$newSheet.Range("G7").Value = '''0.0378'
$newSheet.Range("H7").Value = 9
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = '''005776'
$newSheet.Range("C8").Value = '中加转型动力灵活配置混合C'
$newSheet.Range("D8").Value = '''1.27'
$newSheet.Range("E8").Value = '''50.55'
$newSheet.Range("F8").Value = '''2.87'
$newSheet.Range("G8").Value = '''0.0364'
$newSheet.Range("H8").Value = 5
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = '''012071'
$newSheet.Range("C9").Value = '中加喜利回报一年持有期混合A'
$newSheet.Range("D9").Value = '''1.98'
$newSheet.Range("E9").Value = '''38.64'
$newSheet.Range("F9").Value = '''1.71'
$newSheet.Range("G9").Value = '''0.0339'
$newSheet.Range("H9").Value = 9
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = '''009928'
$newSheet.Range("C10").Value = '工银瑞信聚利18个月定期开放混合C'
$newSheet.Range("D10").Value = '''0.39'
$newSheet.Range("E10").Value = '''26.72'
$newSheet.Range("F10").Value = '''1.94'
$newSheet.Range("G10").Value = '''0.0076'
$newSheet.Range("H10").Value = 3
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = '''851088'
$newSheet.Range("C11").Value = '海通量化成长精选一年持有期混合A'
$newSheet.Range("D11").Value = '''0.38'
$newSheet.Range("E11").Value = '''85.56'
$newSheet.Range("F11").Value = '''1.36'
$newSheet.Range("G11").Value = '''0.0052'
$newSheet.Range("H11").Value = 1
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = '''850010'
$newSheet.Range("C12").Value = '海通量化成长精选一年持有期混合B'
$newSheet.Range("D12").Value = '''0.25'
$newSheet.Range("E12").Value = '''85.56'
$newSheet.Range("F12").Value = '''1.36'
$newSheet.Range("G12").Value = '''0.0034'
$newSheet.Range("H12").Value = 1
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = '''005053'
$newSheet.Range("C13").Value = '银河量化价值混合A'
$newSheet.Range("D13").Value = '''0.10'
$newSheet.Range("E13").Value = '''78.55'
$newSheet.Range("F13").Value = '''1.74'
$newSheet.Range("G13").Value = '''0.0017'
$newSheet.Range("H13").Value = 6
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = '''005126'
$newSheet.Range("C14").Value = '银河量化稳进混合'
$newSheet.Range("D14").Value = '''0.13'
$newSheet.Range("E14").Value = '''55.69'
$newSheet.Range("F14").Value = '''1.14'
$newSheet.Range("G14").Value = '''0.0015'
$newSheet.Range("H14").Value = 4
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = '''851099'
$newSheet.Range("C15").Value = '海通量化成长精选一年持有期混合C'
$newSheet.Range("D15").Value = '''0.03'
$newSheet.Range("E15").Value = '''85.56'
$newSheet.Range("F15").Value = '''1.36'
$newSheet.Range("G15").Value = '''0.0004'
$newSheet.Range("H15").Value = 1
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = '''013026'
$newSheet.Range("C16").Value = '银河量化价值混合C'
$newSheet.Range("D16").Value = '''0.00'
$newSheet.Range("E16").Value = '''78.55'
$newSheet.Range("F16").Value = '''1.74'
$newSheet.Range("G16").Value = 0
$newSheet.Range("H16").Value = 6

# --- 4. Update "总计" (summary) sheet: insert new 2022-Q3 row at the top ---
$ws = $sheetTotal
$ws.Rows.Item(2).Insert()
$ws.Range("B2:D2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 0.52
for ($i = 1; $i -le 7; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
